$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the numeric values in G2:G27 with text labels "Index1_1".."Index1_26"
for ($i = 2; $i -le 27; $i++) {
    $n = $i - 1
    $ws.Cells.Item($i, 7).Value = "Index1_$n"
}

# Update the active selection to reflect the new selection on G2:G27
$ws.Range("G2:G27").Select()
